$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at the last used row (07-09-2021). Append the
# next 8 calendar days (08-09-2021 .. 15-09-2021), repeating the same
# "Monto comprado" values used throughout September.
$lastRow = $ws.UsedRange.Rows.Count
$startRow = $lastRow + 1

$dates = @("08-09-2021", "09-09-2021", "10-09-2021", "11-09-2021", "12-09-2021", "13-09-2021", "14-09-2021", "15-09-2021")

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Force column A to be entered as literal text so Excel doesn't
    # reinterpret a "dd-mm-yyyy"-looking string as a date serial, then
    # drop back to the sheet's normal (unformatted) style so the new
    # cells match the look of the existing rows above them.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = 3623
    $ws.Cells.Item($r, 3).Value = 240
}
